$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.208.43"
$ws.Range("E2").Value = "'  +1.11%  "
$ws.Range("D3").Value = "'3.588.40"
$ws.Range("E3").Value = "'  +0.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("D5").Value = "'581.01"
$ws.Range("E5").Value = "'  -1.44%  "
$ws.Range("D6").Value = "'191.49"
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("D7").Value = "'0.635"
$ws.Range("E7").Value = "'  -1.46%  "
$ws.Range("D8").Value = "'3.579.62"
$ws.Range("E8").Value = "'  +0.40%  "
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("E10").Value = "'  -2.29%  "
$ws.Range("D11").Value = "'0.664"
$ws.Range("E11").Value = "'  +0.46%  "
$ws.Range("D12").Value = "'56.96"
$ws.Range("E12").Value = "'  -1.72%  "
$ws.Range("D13").Value = "'0.0000302"
$ws.Range("E13").Value = "'  +2.05%  "
$ws.Range("D14").Value = "'9.88"
$ws.Range("E14").Value = "'  +2.19%  "
$ws.Range("D15").Value = "'4.165.70"
$ws.Range("E15").Value = "'  +1.19%  "
$ws.Range("D16").Value = "'20.29"
$ws.Range("E16").Value = "'  +5.46%  "
$ws.Range("D17").Value = "'3.586.65"
$ws.Range("E17").Value = "'  +0.76%  "
$ws.Range("D18").Value = "'70.175.56"
$ws.Range("E18").Value = "'  +1.21%  "
$ws.Range("E19").Value = "'  +1.00%  "
$ws.Range("D20").Value = "'0.122"
$ws.Range("E20").Value = "'  +1.46%  "
$ws.Range("D21").Value = "'1.05"
$ws.Range("E21").Value = "'  +0.09%  "
$ws.Range("B22").Value = "'InternetComputer(DFINITY)"
$ws.Range("C22").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'20.23"
$ws.Range("E22").Value = "'  +16.14%  "
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'481.04"
$ws.Range("E23").Value = "'  -2.81%  "
$ws.Range("D24").Value = "'5.14"
$ws.Range("E24").Value = "'  -8.26%  "
$ws.Range("D25").Value = "'4.39"
$ws.Range("E25").Value = "'  -1.06%  "
$ws.Range("D26").Value = "'89.18"
$ws.Range("E26").Value = "'  -2.34%  "
$ws.Range("E27").Value = "'  +1.82%  "
$ws.Range("D28").Value = "'11.28"
$ws.Range("E28").Value = "'  +0.87%  "
$ws.Range("E29").Value = "'  +0.48%  "
$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.80"
$ws.Range("E30").Value = "'  +4.56%  "
$ws.Range("B31").Value = "'EthereumClassic"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'32.29"
$ws.Range("E31").Value = "'  +0.92%  "
$ws.Range("E32").Value = "'  +6.98%  "
$ws.Range("D33").Value = "'66.97"
$ws.Range("E33").Value = "'  +2.56%  "
$ws.Range("D34").Value = "'12.21"
$ws.Range("E34").Value = "'  +0.84%  "
$ws.Range("D35").Value = "'603.80"
$ws.Range("E35").Value = "'  -1.42%  "
$ws.Range("D36").Value = "'40.53"
$ws.Range("E36").Value = "'  +7.26%  "
$ws.Range("D37").Value = "'0.0₃0812"
$ws.Range("E37").Value = "'  -2.29%  "
$ws.Range("D38").Value = "'0.405"
$ws.Range("E38").Value = "'  +2.44%  "
$ws.Range("B39").Value = "'Kaspa"
$ws.Range("C39").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.146"
$ws.Range("E39").Value = "'  -1.41%  "
$ws.Range("B40").Value = "'Dai"
$ws.Range("C40").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "'  -0.25%  "
$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = "'  +12.11%  "
$ws.Range("D42").Value = "'3.56"
$ws.Range("E42").Value = "'  -0.93%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'3.266.33"
$ws.Range("E43").Value = "'  -1.07%  "
$ws.Range("B44").Value = "'dogwifhat"
$ws.Range("C44").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = "'  +12.45%  "
$ws.Range("D45").Value = "'3.15"
$ws.Range("E45").Value = "'  +3.26%  "
$ws.Range("E46").Value = "'  +1.89%  "
$ws.Range("D47").Value = "'9.63"
$ws.Range("E47").Value = "'  +6.49%  "
$ws.Range("E48").Value = "'  +2.28%  "
$ws.Range("E49").Value = "'  +0.16%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "'  +0.08%  "
$ws.Range("E51").Value = "'  -0.52%  "
